# Insert a new daily price record for "Feria Lagunitas de Puerto Montt - Cilantro"
# as a new row 278, pushing the existing rows 278-351 down to 279-352.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A278").EntireRow.Insert()

$ws.Range("A278").Value = 4
$ws.Range("B278").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C278").Value = "Los Lagos"
$ws.Range("D278").Value = 44855
$ws.Range("E278").Value = 10
$ws.Range("F278").Value = 100112040
$ws.Range("G278").Value = "Cilantro"
$ws.Range("H278").Value = "Sin especificar"
$ws.Range("I278").Value = "Primera"
$ws.Range("J278").Value = 240
$ws.Range("K278").Value = 10000
$ws.Range("L278").Value = 10500
$ws.Range("M278").Value = 10250
$ws.Range("N278").Value = "$/caja 36 atados"
$ws.Range("O278").Value = "Región Metropolitana"
$ws.Range("P278").Value = 285
$ws.Range("Q278").Value = 36
$ws.Range("R278").Value = "Hortaliza"
